$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.233.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.555.77'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.07'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.90'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.554.45'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.487'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.33%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.85'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.409'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.161.61'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.24'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.549.71'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.118'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.269.97'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.05'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.26'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.66%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '419.66'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.606'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.85'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.697.61'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.23'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.92'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.48'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.553.33'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.66'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.19%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.65'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.50%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.33'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -8.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.36'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.63%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '174.27'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0830'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.13'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.865'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.83%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.83'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.45'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.18'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.12'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.63'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.00%  '
